$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 162-172 with new/shifted weekly price data,
# and append new rows 173-175 (duplicated week pushed down).

# Row 162
$ws.Cells.Item(162, 1).Value = 11
$ws.Cells.Item(162, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(162, 3).Value = 'Bíobío'
$ws.Cells.Item(162, 4).Value = 44491
$ws.Cells.Item(162, 5).Value = 8
$ws.Cells.Item(162, 6).Value = 'Fruta'
$ws.Cells.Item(162, 7).Value = 100101
$ws.Cells.Item(162, 8).Value = 'Berries'
$ws.Cells.Item(162, 9).Value = 100112025
$ws.Cells.Item(162, 10).Value = 'Frutilla'
$ws.Cells.Item(162, 11).Value = 'Sin especificar'
$ws.Cells.Item(162, 12).Value = 'Especial'
$ws.Cells.Item(162, 13).Value = 270
$ws.Cells.Item(162, 14).Value = 9000
$ws.Cells.Item(162, 15).Value = 10000
$ws.Cells.Item(162, 16).Value = 9556
$ws.Cells.Item(162, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(162, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(162, 19).Value = 1365
$ws.Cells.Item(162, 20).Value = 7

# Row 163
$ws.Cells.Item(163, 1).Value = 11
$ws.Cells.Item(163, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(163, 3).Value = 'Bíobío'
$ws.Cells.Item(163, 4).Value = 44491
$ws.Cells.Item(163, 5).Value = 8
$ws.Cells.Item(163, 6).Value = 'Fruta'
$ws.Cells.Item(163, 7).Value = 100101
$ws.Cells.Item(163, 8).Value = 'Berries'
$ws.Cells.Item(163, 9).Value = 100112025
$ws.Cells.Item(163, 10).Value = 'Frutilla'
$ws.Cells.Item(163, 11).Value = 'Sin especificar'
$ws.Cells.Item(163, 12).Value = 'Primera'
$ws.Cells.Item(163, 13).Value = 380
$ws.Cells.Item(163, 14).Value = 6500
$ws.Cells.Item(163, 15).Value = 7000
$ws.Cells.Item(163, 16).Value = 6737
$ws.Cells.Item(163, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(163, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(163, 19).Value = 962
$ws.Cells.Item(163, 20).Value = 7

# Row 164
$ws.Cells.Item(164, 1).Value = 11
$ws.Cells.Item(164, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(164, 3).Value = 'Bíobío'
$ws.Cells.Item(164, 4).Value = 44491
$ws.Cells.Item(164, 5).Value = 8
$ws.Cells.Item(164, 6).Value = 'Fruta'
$ws.Cells.Item(164, 7).Value = 100101
$ws.Cells.Item(164, 8).Value = 'Berries'
$ws.Cells.Item(164, 9).Value = 100112025
$ws.Cells.Item(164, 10).Value = 'Frutilla'
$ws.Cells.Item(164, 11).Value = 'Sin especificar'
$ws.Cells.Item(164, 12).Value = 'Segunda'
$ws.Cells.Item(164, 13).Value = 200
$ws.Cells.Item(164, 14).Value = 5500
$ws.Cells.Item(164, 15).Value = 5500
$ws.Cells.Item(164, 16).Value = 5500
$ws.Cells.Item(164, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(164, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(164, 19).Value = 786
$ws.Cells.Item(164, 20).Value = 7

# Row 165
$ws.Cells.Item(165, 1).Value = 11
$ws.Cells.Item(165, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(165, 3).Value = 'Bíobío'
$ws.Cells.Item(165, 4).Value = 44328
$ws.Cells.Item(165, 5).Value = 8
$ws.Cells.Item(165, 6).Value = 'Fruta'
$ws.Cells.Item(165, 7).Value = 100101
$ws.Cells.Item(165, 8).Value = 'Berries'
$ws.Cells.Item(165, 9).Value = 100112025
$ws.Cells.Item(165, 10).Value = 'Frutilla'
$ws.Cells.Item(165, 11).Value = 'Sin especificar'
$ws.Cells.Item(165, 12).Value = 'Especial'
$ws.Cells.Item(165, 13).Value = 100
$ws.Cells.Item(165, 14).Value = 15000
$ws.Cells.Item(165, 15).Value = 15000
$ws.Cells.Item(165, 16).Value = 15000
$ws.Cells.Item(165, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(165, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(165, 19).Value = 2143
$ws.Cells.Item(165, 20).Value = 7

# Row 166
$ws.Cells.Item(166, 1).Value = 11
$ws.Cells.Item(166, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(166, 3).Value = 'Bíobío'
$ws.Cells.Item(166, 4).Value = 44328
$ws.Cells.Item(166, 5).Value = 8
$ws.Cells.Item(166, 6).Value = 'Fruta'
$ws.Cells.Item(166, 7).Value = 100101
$ws.Cells.Item(166, 8).Value = 'Berries'
$ws.Cells.Item(166, 9).Value = 100112025
$ws.Cells.Item(166, 10).Value = 'Frutilla'
$ws.Cells.Item(166, 11).Value = 'Sin especificar'
$ws.Cells.Item(166, 12).Value = 'Primera'
$ws.Cells.Item(166, 13).Value = 100
$ws.Cells.Item(166, 14).Value = 12000
$ws.Cells.Item(166, 15).Value = 12000
$ws.Cells.Item(166, 16).Value = 12000
$ws.Cells.Item(166, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(166, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(166, 19).Value = 1714
$ws.Cells.Item(166, 20).Value = 7

# Row 167
$ws.Cells.Item(167, 1).Value = 11
$ws.Cells.Item(167, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(167, 3).Value = 'Bíobío'
$ws.Cells.Item(167, 4).Value = 44223
$ws.Cells.Item(167, 5).Value = 8
$ws.Cells.Item(167, 6).Value = 'Fruta'
$ws.Cells.Item(167, 7).Value = 100101
$ws.Cells.Item(167, 8).Value = 'Berries'
$ws.Cells.Item(167, 9).Value = 100112025
$ws.Cells.Item(167, 10).Value = 'Frutilla'
$ws.Cells.Item(167, 11).Value = 'Sin especificar'
$ws.Cells.Item(167, 12).Value = 'Especial'
$ws.Cells.Item(167, 13).Value = 200
$ws.Cells.Item(167, 14).Value = 10000
$ws.Cells.Item(167, 15).Value = 10000
$ws.Cells.Item(167, 16).Value = 10000
$ws.Cells.Item(167, 17).Value = '$/bandeja 3 kilos'
$ws.Cells.Item(167, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(167, 19).Value = 3333
$ws.Cells.Item(167, 20).Value = 3

# Row 168
$ws.Cells.Item(168, 1).Value = 11
$ws.Cells.Item(168, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(168, 3).Value = 'Bíobío'
$ws.Cells.Item(168, 4).Value = 44223
$ws.Cells.Item(168, 5).Value = 8
$ws.Cells.Item(168, 6).Value = 'Fruta'
$ws.Cells.Item(168, 7).Value = 100101
$ws.Cells.Item(168, 8).Value = 'Berries'
$ws.Cells.Item(168, 9).Value = 100112025
$ws.Cells.Item(168, 10).Value = 'Frutilla'
$ws.Cells.Item(168, 11).Value = 'Sin especificar'
$ws.Cells.Item(168, 12).Value = 'Primera'
$ws.Cells.Item(168, 13).Value = 200
$ws.Cells.Item(168, 14).Value = 8000
$ws.Cells.Item(168, 15).Value = 8000
$ws.Cells.Item(168, 16).Value = 8000
$ws.Cells.Item(168, 17).Value = '$/bandeja 3 kilos'
$ws.Cells.Item(168, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(168, 19).Value = 2667
$ws.Cells.Item(168, 20).Value = 3

# Row 169
$ws.Cells.Item(169, 1).Value = 11
$ws.Cells.Item(169, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(169, 3).Value = 'Bíobío'
$ws.Cells.Item(169, 4).Value = 44223
$ws.Cells.Item(169, 5).Value = 8
$ws.Cells.Item(169, 6).Value = 'Fruta'
$ws.Cells.Item(169, 7).Value = 100101
$ws.Cells.Item(169, 8).Value = 'Berries'
$ws.Cells.Item(169, 9).Value = 100112025
$ws.Cells.Item(169, 10).Value = 'Frutilla'
$ws.Cells.Item(169, 11).Value = 'Sin especificar'
$ws.Cells.Item(169, 12).Value = 'Segunda'
$ws.Cells.Item(169, 13).Value = 100
$ws.Cells.Item(169, 14).Value = 6000
$ws.Cells.Item(169, 15).Value = 6000
$ws.Cells.Item(169, 16).Value = 6000
$ws.Cells.Item(169, 17).Value = '$/bandeja 3 kilos'
$ws.Cells.Item(169, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(169, 19).Value = 2000
$ws.Cells.Item(169, 20).Value = 3

# Row 170
$ws.Cells.Item(170, 1).Value = 11
$ws.Cells.Item(170, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(170, 3).Value = 'Bíobío'
$ws.Cells.Item(170, 4).Value = 44187
$ws.Cells.Item(170, 5).Value = 8
$ws.Cells.Item(170, 6).Value = 'Fruta'
$ws.Cells.Item(170, 7).Value = 100101
$ws.Cells.Item(170, 8).Value = 'Berries'
$ws.Cells.Item(170, 9).Value = 100112025
$ws.Cells.Item(170, 10).Value = 'Frutilla'
$ws.Cells.Item(170, 11).Value = 'Sin especificar'
$ws.Cells.Item(170, 12).Value = 'Especial'
$ws.Cells.Item(170, 13).Value = 100
$ws.Cells.Item(170, 14).Value = 11000
$ws.Cells.Item(170, 15).Value = 11000
$ws.Cells.Item(170, 16).Value = 11000
$ws.Cells.Item(170, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(170, 18).Value = 'Región del Maule'
$ws.Cells.Item(170, 19).Value = 1571
$ws.Cells.Item(170, 20).Value = 7

# Row 171
$ws.Cells.Item(171, 1).Value = 11
$ws.Cells.Item(171, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(171, 3).Value = 'Bíobío'
$ws.Cells.Item(171, 4).Value = 44187
$ws.Cells.Item(171, 5).Value = 8
$ws.Cells.Item(171, 6).Value = 'Fruta'
$ws.Cells.Item(171, 7).Value = 100101
$ws.Cells.Item(171, 8).Value = 'Berries'
$ws.Cells.Item(171, 9).Value = 100112025
$ws.Cells.Item(171, 10).Value = 'Frutilla'
$ws.Cells.Item(171, 11).Value = 'Sin especificar'
$ws.Cells.Item(171, 12).Value = 'Primera'
$ws.Cells.Item(171, 13).Value = 100
$ws.Cells.Item(171, 14).Value = 9000
$ws.Cells.Item(171, 15).Value = 9000
$ws.Cells.Item(171, 16).Value = 9000
$ws.Cells.Item(171, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(171, 18).Value = 'Región del Maule'
$ws.Cells.Item(171, 19).Value = 1286
$ws.Cells.Item(171, 20).Value = 7

# Row 172
$ws.Cells.Item(172, 1).Value = 11
$ws.Cells.Item(172, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(172, 3).Value = 'Bíobío'
$ws.Cells.Item(172, 4).Value = 44187
$ws.Cells.Item(172, 5).Value = 8
$ws.Cells.Item(172, 6).Value = 'Fruta'
$ws.Cells.Item(172, 7).Value = 100101
$ws.Cells.Item(172, 8).Value = 'Berries'
$ws.Cells.Item(172, 9).Value = 100112025
$ws.Cells.Item(172, 10).Value = 'Frutilla'
$ws.Cells.Item(172, 11).Value = 'Sin especificar'
$ws.Cells.Item(172, 12).Value = 'Segunda'
$ws.Cells.Item(172, 13).Value = 50
$ws.Cells.Item(172, 14).Value = 7000
$ws.Cells.Item(172, 15).Value = 7000
$ws.Cells.Item(172, 16).Value = 7000
$ws.Cells.Item(172, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(172, 18).Value = 'Región del Maule'
$ws.Cells.Item(172, 19).Value = 1000
$ws.Cells.Item(172, 20).Value = 7

# Row 173
$ws.Cells.Item(173, 1).Value = 11
$ws.Cells.Item(173, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(173, 3).Value = 'Bíobío'
$ws.Cells.Item(173, 4).Value = 44250
$ws.Cells.Item(173, 5).Value = 8
$ws.Cells.Item(173, 6).Value = 'Fruta'
$ws.Cells.Item(173, 7).Value = 100101
$ws.Cells.Item(173, 8).Value = 'Berries'
$ws.Cells.Item(173, 9).Value = 100112025
$ws.Cells.Item(173, 10).Value = 'Frutilla'
$ws.Cells.Item(173, 11).Value = 'Sin especificar'
$ws.Cells.Item(173, 12).Value = 'Especial'
$ws.Cells.Item(173, 13).Value = 200
$ws.Cells.Item(173, 14).Value = 10000
$ws.Cells.Item(173, 15).Value = 10000
$ws.Cells.Item(173, 16).Value = 10000
$ws.Cells.Item(173, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(173, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(173, 19).Value = 1429
$ws.Cells.Item(173, 20).Value = 7
$ws.Cells.Item(173, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 174
$ws.Cells.Item(174, 1).Value = 11
$ws.Cells.Item(174, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(174, 3).Value = 'Bíobío'
$ws.Cells.Item(174, 4).Value = 44250
$ws.Cells.Item(174, 5).Value = 8
$ws.Cells.Item(174, 6).Value = 'Fruta'
$ws.Cells.Item(174, 7).Value = 100101
$ws.Cells.Item(174, 8).Value = 'Berries'
$ws.Cells.Item(174, 9).Value = 100112025
$ws.Cells.Item(174, 10).Value = 'Frutilla'
$ws.Cells.Item(174, 11).Value = 'Sin especificar'
$ws.Cells.Item(174, 12).Value = 'Primera'
$ws.Cells.Item(174, 13).Value = 200
$ws.Cells.Item(174, 14).Value = 8000
$ws.Cells.Item(174, 15).Value = 8000
$ws.Cells.Item(174, 16).Value = 8000
$ws.Cells.Item(174, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(174, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(174, 19).Value = 1143
$ws.Cells.Item(174, 20).Value = 7
$ws.Cells.Item(174, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 175
$ws.Cells.Item(175, 1).Value = 11
$ws.Cells.Item(175, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(175, 3).Value = 'Bíobío'
$ws.Cells.Item(175, 4).Value = 44250
$ws.Cells.Item(175, 5).Value = 8
$ws.Cells.Item(175, 6).Value = 'Fruta'
$ws.Cells.Item(175, 7).Value = 100101
$ws.Cells.Item(175, 8).Value = 'Berries'
$ws.Cells.Item(175, 9).Value = 100112025
$ws.Cells.Item(175, 10).Value = 'Frutilla'
$ws.Cells.Item(175, 11).Value = 'Sin especificar'
$ws.Cells.Item(175, 12).Value = 'Segunda'
$ws.Cells.Item(175, 13).Value = 50
$ws.Cells.Item(175, 14).Value = 6000
$ws.Cells.Item(175, 15).Value = 6000
$ws.Cells.Item(175, 16).Value = 6000
$ws.Cells.Item(175, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(175, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(175, 19).Value = 857
$ws.Cells.Item(175, 20).Value = 7
$ws.Cells.Item(175, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
